$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J12").Value = "2020-007"
